$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cerinte Generale")

# Update status / responsible columns on row 11
$ws.Range("E11").Value = "Alex"
$ws.Range("F11").Value = "in lucru"

# Update the sheet view: scroll back to top-left (A1) and select F12
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F12").Select()
